$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in salary / salary-year (and one notes) for ranks 41-50 (rows 42-51) ---
# Row 42 - Colorado Springs, Colorado
$ws.Range("J42").Value = 114159
$ws.Range("K42").Value = "FY2019"

# Row 43 - Omaha, Nebraska
$ws.Range("J43").Value = 104358
$ws.Range("K43").Value = "FY2019"
$ws.Range("M43:N43").NumberFormat = $ws.Range("J2").NumberFormat

# Row 44 - Raleigh, North Carolina
$ws.Range("J44").Value = 23720
$ws.Range("K44").Value = "FY2019"

# Row 45 - Miami, Florida
$ws.Range("J45").Value = 97000
$ws.Range("K45").Value = "FY2019"

# Row 46 - Oakland, California
$ws.Range("J46").Value = 212000
$ws.Range("K46").Value = "FY2018"

# Row 47 - Minneapolis, Minnesota
$ws.Range("J47").Value = 126528
$ws.Range("K47").Value = "FY2018"

# Row 48 - Tulsa, Oklahoma
$ws.Range("J48").Value = 105000
$ws.Range("K48").Value = "FY2018"

# Row 49 - Cleveland, Ohio
$ws.Range("J49").Value = 140888.56
$ws.Range("K49").Value = "FY2017"

# Row 50 - Wichita, Kansas
$ws.Range("J50").Value = 103560
$ws.Range("K50").Value = "FY2019"

# Row 51 - Arlington, Texas (no official salary-year source; note instead)
$ws.Range("J51").Value = 3000
$ws.Range("L51").Value = "In city charter as $250 per month"

# --- Row 66 - St. Paul, Minnesota also gets a salary added ---
$ws.Range("J66").Value = 129000
$ws.Range("K66").Value = "FY2019"

# --- New columns M:N (13:14) sized to match column J ---
$ws.Columns("M:N").ColumnWidth = 11.67

# --- Update the selection / view to where editing left off ---
[void]$ws.Range("J52").Select()
$excel.ActiveWindow.Zoom = 100
